# Fix marksheet so it correctly handles (float) marks input instead of
# leaving stale placeholder figures / text-typed numbers, and rebuild the
# per-question "Student Ans / Correct Ans" table as a single two-column
# block (A:B) instead of the old broken 3-block (A:B, D:E, G:H) layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Score summary block (rows 10-12)
# ---------------------------------------------------------------------

# Row 10 - "No." (counts)
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value2 = 25
$ws.Range("C10").Value2 = 1
$ws.Range("D10").Value2 = 2
$ws.Range("E10").Value2 = 28

# Row 11 - "Marking" (per-question weight); C11 used to be stored as the
# text "-1" (a leftover string) instead of a real number - make it numeric.
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value2 = 4
$ws.Range("C11").Value2 = -1

# Row 12 - "Total"
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value2 = 100
$ws.Range("C12").Value2 = -1
$ws.Range("E12").Value2 = "99/112"

# ---------------------------------------------------------------------
# 2) Rebuild the answers table
# ---------------------------------------------------------------------
# The "Correct Ans" column (B) already holds the right data for all 25
# questions (rows 16-40) and is left untouched. The "Student Ans" column
# (A) used to only be filled for a handful of rows, with the remaining
# answers wrongly duplicated into extra D:E / G:H column blocks. Collapse
# everything back into column A, colouring each answer green
# (correctStyle) when it matches the correct answer, red (incorrectStyle)
# when it does not, and leaving it blank (normalStyle) when the student
# did not attempt the question.

$correctAnswers = @{
    16 = "Option A"; 17 = "Option D"; 18 = "Option B"; 19 = "Option C"; 20 = "Option B";
    21 = "Option C"; 22 = "Option D"; 23 = "Option D"; 24 = "Option A"; 25 = "Option A";
    26 = "Option C"; 27 = "Option A"; 28 = "Option D"; 29 = "Option D"; 30 = "Option B";
    31 = "Option D"; 32 = "Option C"; 33 = "Option D"; 34 = "Option B"; 35 = "Option D";
    36 = "Option A"; 37 = "Option A"; 38 = "Option A"; 39 = "Option D"; 40 = "Option D";
}

# Student answers; $null means the question was not attempted.
$studentAnswers = @{
    16 = "Option A"; 17 = "Option D"; 18 = "Option B"; 19 = "Option C"; 20 = "Option B";
    21 = "Option C"; 22 = "Option D"; 23 = $null;       24 = $null;       25 = "Option A";
    26 = "Option C"; 27 = "Option A"; 28 = "Option D"; 29 = "Option D"; 30 = "Option B";
    31 = "Option D"; 32 = "Option C"; 33 = "Option D"; 34 = "Option A"; 35 = "Option D";
    36 = "Option A"; 37 = "Option A"; 38 = "Option A"; 39 = "Option D"; 40 = "Option D";
}

for ($row = 16; $row -le 40; $row++) {
    $cell = $ws.Range("A$row")
    $student = $studentAnswers[$row]
    $correct = $correctAnswers[$row]

    if ($null -eq $student) {
        $cell.Value2 = ""
        $cell.Style = "normalStyle"
    } elseif ($student -eq $correct) {
        $cell.Value2 = $student
        $cell.Style = "correctStyle"
    } else {
        $cell.Value2 = $student
        $cell.Style = "incorrectStyle"
    }
}

# The first three questions of the table used to additionally live in the
# D:E block - keep that duplicate in sync with the (now single) table
# instead of the independent/stale values it used to show.
$ws.Range("D16").Value2 = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D17").Value2 = "Option C"
$ws.Range("D17").Style = "correctStyle"
$ws.Range("D18").Value2 = "Option D"
$ws.Range("D18").Style = "correctStyle"

# Drop the rest of the old duplicated D:E block and the whole G:H block -
# everything now lives in A:B.
$ws.Range("D19:E40").Clear()
$ws.Range("G15:H40").Clear()

